$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''30.274.72'
$ws.Range("D2").Style = 'Normal'
$ws.Range("E2").Value = '''  +1.40%  '
$ws.Range("E2").Style = 'Normal'

$ws.Range("D3").Value = '''1.897.40'
$ws.Range("D3").Style = 'Normal'
$ws.Range("E3").Value = '''  +0.53%  '
$ws.Range("E3").Style = 'Normal'

$ws.Range("D4").Value = '''0.9986'
$ws.Range("D4").Style = 'Normal'
$ws.Range("E4").Value = '''  -0.17%  '
$ws.Range("E4").Style = 'Normal'

$ws.Range("D5").Value = '''0.7491'
$ws.Range("D5").Style = 'Normal'
$ws.Range("E5").Value = '''  +0.15%  '
$ws.Range("E5").Style = 'Normal'

$ws.Range("D6").Value = '''243.10'
$ws.Range("D6").Style = 'Normal'
$ws.Range("E6").Value = '''  +0.24%  '
$ws.Range("E6").Style = 'Normal'

$ws.Range("D7").Value = '''0.9991'
$ws.Range("D7").Style = 'Normal'
$ws.Range("E7").Value = '''  -0.15%  '
$ws.Range("E7").Style = 'Normal'

$ws.Range("D8").Value = '''0.3181'
$ws.Range("D8").Style = 'Normal'
$ws.Range("E8").Value = '''  +1.98%  '
$ws.Range("E8").Style = 'Normal'

$ws.Range("D9").Value = '''0.07267'
$ws.Range("D9").Style = 'Normal'
$ws.Range("E9").Value = '''  +1.89%  '
$ws.Range("E9").Style = 'Normal'

$ws.Range("D10").Value = '''25.14'
$ws.Range("D10").Style = 'Normal'
$ws.Range("E10").Value = '''  -1.20%  '
$ws.Range("E10").Style = 'Normal'

$ws.Range("D11").Value = '''0.08376'
$ws.Range("D11").Style = 'Normal'
$ws.Range("E11").Value = '''  -1.20%  '
$ws.Range("E11").Style = 'Normal'

$ws.Range("D12").Value = '''0.7649'
$ws.Range("D12").Style = 'Normal'
$ws.Range("E12").Value = '''  +0.61%  '
$ws.Range("E12").Style = 'Normal'

$ws.Range("B13").Value = '''WrappedEther'
$ws.Range("B13").Style = 'Normal'
$ws.Range("C13").Value = '''https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("C13").Style = 'Normal'
$ws.Range("D13").Value = '''1.930.08'
$ws.Range("D13").Style = 'Normal'
$ws.Range("E13").Value = '''  +0.62%  '
$ws.Range("E13").Style = 'Normal'

$ws.Range("B14").Value = '''Polkadot'
$ws.Range("B14").Style = 'Normal'
$ws.Range("C14").Value = '''https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("C14").Style = 'Normal'
$ws.Range("D14").Value = '''5.470'
$ws.Range("D14").Style = 'Normal'
$ws.Range("E14").Value = '''  +2.03%  '
$ws.Range("E14").Style = 'Normal'

$ws.Range("D15").Value = '''93.33'
$ws.Range("D15").Style = 'Normal'
$ws.Range("E15").Value = '''  -0.08%  '
$ws.Range("E15").Style = 'Normal'

$ws.Range("D16").Value = '''6.207'
$ws.Range("D16").Style = 'Normal'
$ws.Range("E16").Value = '''  +0.97%  '
$ws.Range("E16").Style = 'Normal'

$ws.Range("D17").Value = '''30.274.13'
$ws.Range("D17").Style = 'Normal'
$ws.Range("E17").Value = '''  +1.08%  '
$ws.Range("E17").Style = 'Normal'

$ws.Range("D18").Value = '''252.12'
$ws.Range("D18").Style = 'Normal'

$ws.Range("E19").Value = '''  +0.02%  '
$ws.Range("E19").Style = 'Normal'

$ws.Range("D20").Value = '''0.000007893'
$ws.Range("D20").Style = 'Normal'
$ws.Range("E20").Value = '''  +1.31%  '
$ws.Range("E20").Style = 'Normal'

$ws.Range("D21").Value = '''2.152.37'
$ws.Range("D21").Style = 'Normal'
$ws.Range("E21").Value = '''  -0.45%  '
$ws.Range("E21").Style = 'Normal'

$ws.Range("D22").Value = '''0.9978'
$ws.Range("D22").Style = 'Normal'
$ws.Range("E22").Value = '''  -0.14%  '
$ws.Range("E22").Style = 'Normal'

$ws.Range("D23").Value = '''8.042'
$ws.Range("D23").Style = 'Normal'
$ws.Range("E23").Value = '''  +0.25%  '
$ws.Range("E23").Style = 'Normal'

$ws.Range("E24").Value = '''  -0.20%  '
$ws.Range("E24").Style = 'Normal'

$ws.Range("D25").Value = '''0.1606'
$ws.Range("D25").Style = 'Normal'
$ws.Range("E25").Value = '''  +0.85%  '
$ws.Range("E25").Style = 'Normal'

$ws.Range("D26").Value = '''9.336'
$ws.Range("D26").Style = 'Normal'
$ws.Range("E26").Value = '''  -0.45%  '
$ws.Range("E26").Style = 'Normal'

$ws.Range("E27").Value = '''  +1.28%  '
$ws.Range("E27").Style = 'Normal'

$ws.Range("D28").Value = '''18.86'
$ws.Range("D28").Style = 'Normal'
$ws.Range("E28").Value = '''  +0.53%  '
$ws.Range("E28").Style = 'Normal'

$ws.Range("D29").Value = '''2.084'
$ws.Range("D29").Style = 'Normal'
$ws.Range("E29").Value = '''  +2.69%  '
$ws.Range("E29").Style = 'Normal'

$ws.Range("D30").Value = '''1.478'
$ws.Range("D30").Style = 'Normal'
$ws.Range("E30").Value = '''  -2.34%  '
$ws.Range("E30").Style = 'Normal'

$ws.Range("D31").Value = '''4.620'
$ws.Range("D31").Style = 'Normal'
$ws.Range("E31").Value = '''  +3.20%  '
$ws.Range("E31").Style = 'Normal'

$ws.Range("D32").Value = '''1.542'
$ws.Range("D32").Style = 'Normal'
$ws.Range("E32").Value = '''  +0.76%  '
$ws.Range("E32").Style = 'Normal'

$ws.Range("D33").Value = '''4.241'
$ws.Range("D33").Style = 'Normal'
$ws.Range("E33").Value = '''  +3.31%  '
$ws.Range("E33").Style = 'Normal'

$ws.Range("D34").Value = '''0.05446'
$ws.Range("D34").Style = 'Normal'
$ws.Range("E34").Value = '''  +0.94%  '
$ws.Range("E34").Style = 'Normal'

$ws.Range("D35").Value = '''1.266'
$ws.Range("D35").Style = 'Normal'
$ws.Range("E35").Value = '''  +2.34%  '
$ws.Range("E35").Style = 'Normal'

$ws.Range("D36").Value = '''0.7680'
$ws.Range("D36").Style = 'Normal'
$ws.Range("E36").Value = '''  +3.23%  '
$ws.Range("E36").Style = 'Normal'

$ws.Range("D37").Value = '''0.9972'
$ws.Range("D37").Style = 'Normal'
$ws.Range("E37").Value = '''  -0.72%  '
$ws.Range("E37").Style = 'Normal'

$ws.Range("D38").Value = '''2.718'
$ws.Range("D38").Style = 'Normal'
$ws.Range("E38").Value = '''  +0.29%  '
$ws.Range("E38").Style = 'Normal'

$ws.Range("D39").Value = '''0.01984'
$ws.Range("D39").Style = 'Normal'
$ws.Range("E39").Value = '''  +2.66%  '
$ws.Range("E39").Style = 'Normal'

$ws.Range("E40").Value = '''  +0.20%  '
$ws.Range("E40").Style = 'Normal'

$ws.Range("D41").Value = '''0.4601'
$ws.Range("D41").Style = 'Normal'
$ws.Range("E41").Value = '''  +3.31%  '
$ws.Range("E41").Style = 'Normal'

$ws.Range("D42").Value = '''1.102.41'
$ws.Range("D42").Style = 'Normal'
$ws.Range("E42").Value = '''  +0.62%  '
$ws.Range("E42").Style = 'Normal'

$ws.Range("D43").Value = '''6.101'
$ws.Range("D43").Style = 'Normal'
$ws.Range("E43").Value = '''  +0.39%  '
$ws.Range("E43").Style = 'Normal'

$ws.Range("D44").Value = '''73.28'
$ws.Range("D44").Style = 'Normal'
$ws.Range("E44").Value = '''  +0.89%  '
$ws.Range("E44").Style = 'Normal'

$ws.Range("D45").Value = '''0.8738'
$ws.Range("D45").Style = 'Normal'
$ws.Range("E45").Value = '''  +1.62%  '
$ws.Range("E45").Style = 'Normal'

$ws.Range("D46").Value = '''104.74'
$ws.Range("D46").Style = 'Normal'
$ws.Range("E46").Value = '''  +2.14%  '
$ws.Range("E46").Style = 'Normal'

$ws.Range("D47").Value = '''1.000'
$ws.Range("D47").Style = 'Normal'
$ws.Range("E47").Value = '''  +0.05%  '
$ws.Range("E47").Style = 'Normal'

$ws.Range("D48").Value = '''1.881'
$ws.Range("D48").Style = 'Normal'
$ws.Range("E48").Value = '''  +1.02%  '
$ws.Range("E48").Style = 'Normal'

$ws.Range("D49").Value = '''7.656'
$ws.Range("D49").Style = 'Normal'
$ws.Range("E49").Value = '''  -0.20%  '
$ws.Range("E49").Style = 'Normal'

$ws.Range("D50").Value = '''9.681'
$ws.Range("D50").Style = 'Normal'
$ws.Range("E50").Value = '''  -0.51%  '
$ws.Range("E50").Style = 'Normal'

$ws.Range("D51").Value = '''2.055.26'
$ws.Range("D51").Style = 'Normal'
$ws.Range("E51").Value = '''  -0.16%  '
$ws.Range("E51").Style = 'Normal'
